$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.127.28"
$ws.Range("E2").Value = "  +4.46%  "
$ws.Range("D3").Value = "3.269.28"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.97"
$ws.Range("E5").Value = "  +3.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "177.45"
$ws.Range("E6").Value = "  +3.90%  "
$ws.Range("E7").Value = "  -0.38%  "
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "3.271.23"
$ws.Range("E9").Value = "  +3.13%  "
$ws.Range("E10").Value = "  +4.58%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.415"
$ws.Range("E12").Value = "  +4.94%  "
$ws.Range("D13").Value = "3.838.21"
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("E14").Value = "  +0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.20"
$ws.Range("E15").Value = "  +3.04%  "
$ws.Range("D16").Value = "67.157.45"
$ws.Range("E16").Value = "  +4.47%  "
$ws.Range("E17").Value = "  +3.42%  "
$ws.Range("D18").Value = "3.271.29"
$ws.Range("E18").Value = "  +3.15%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.86"
$ws.Range("E19").Value = "  +2.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.48"
$ws.Range("E20").Value = "  +3.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "372.89"
$ws.Range("E21").Value = "  +5.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.64"
$ws.Range("E22").Value = "  +6.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.66"
$ws.Range("E23").Value = "  +3.60%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.513"
$ws.Range("E25").Value = "  +2.02%  "
$ws.Range("D26").Value = "3.411.19"
$ws.Range("E26").Value = "  +3.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000119"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.87"
$ws.Range("E28").Value = "  +2.60%  "
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "5.65"
$ws.Range("E32").Value = "  +0.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.68"
$ws.Range("E33").Value = "  +2.51%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  +6.38%  "
$ws.Range("E36").Value = "  +2.93%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "166.99"
$ws.Range("E37").Value = "  +8.47%  "
$ws.Range("E38").Value = "  +4.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.860"
$ws.Range("E39").Value = "  +6.20%  "
$ws.Range("E40").Value = "  +11.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.45"
$ws.Range("E41").Value = "  +5.33%  "
$ws.Range("E42").Value = "  +2.20%  "
$ws.Range("D43").Value = "2.748.69"
$ws.Range("E43").Value = "  +4.48%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.47"
$ws.Range("E44").Value = "  +7.58%  "
$ws.Range("E45").Value = "  +5.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "349.33"
$ws.Range("E46").Value = "  +4.49%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "25.11"
$ws.Range("E47").Value = "  +5.15%  "
$ws.Range("B48").Value = "Hedera"
$ws.Range("C48").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0679"
$ws.Range("E48").Value = "  +3.47%  "
$ws.Range("B49").Value = "OKB"
$ws.Range("C49").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "40.56"
$ws.Range("E49").Value = "  +4.31%  "
$ws.Range("E50").Value = "  +4.11%  "
$ws.Range("E51").Value = "  +2.58%  "
